$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.309.05"
$ws.Range("E2").Value = "  -0.38%  "

$ws.Range("D3").Value = "1.882.99"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.94"
$ws.Range("E5").Value = "  -0.58%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("E7").Value = "  -1.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2818"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06555"
$ws.Range("E9").Value = "  -1.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.73"
$ws.Range("E10").Value = "  +5.36%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07770"
$ws.Range("E11").Value = "  +0.68%  "

$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.79"
$ws.Range("E12").Value = "  -2.19%  "

$ws.Range("D13").Value = "1.891.92"
$ws.Range("E13").Value = "  -0.77%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.130"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6646"
$ws.Range("E15").Value = "  -0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "282.65"
$ws.Range("E16").Value = "  +10.86%  "

$ws.Range("D17").Value = "30.326.99"
$ws.Range("E17").Value = "  -0.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "2.136.60"
$ws.Range("E19").Value = "  -0.80%  "

$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.60"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000007297"
$ws.Range("E21").Value = "  -2.02%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.349"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.157"
$ws.Range("E24").Value = "  -2.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.52"
$ws.Range("E25").Value = "  +0.18%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.261"
$ws.Range("E26").Value = "  -1.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.02"
$ws.Range("E27").Value = "  +0.72%  "

$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.371"
$ws.Range("E29").Value = "  -0.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09810"
$ws.Range("E30").Value = "  -2.82%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.461"
$ws.Range("E31").Value = "  -3.63%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.493"
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.166"
$ws.Range("E33").Value = "  -2.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04706"
$ws.Range("E34").Value = "  -0.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7067"
$ws.Range("E35").Value = "  -3.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.091"
$ws.Range("E36").Value = "  -1.72%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.719"
$ws.Range("E37").Value = "  +0.87%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01866"
$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.728"
$ws.Range("E39").Value = "  +7.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.521"
$ws.Range("E40").Value = "  -2.80%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.09"
$ws.Range("E41").Value = "  -1.55%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8717"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.965"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.16"
$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  +0.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4192"
$ws.Range("E46").Value = "  -0.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "989.98"
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.199"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.268"
$ws.Range("E49").Value = "  +5.15%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1161"
$ws.Range("E50").Value = "  -2.84%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.06"
$ws.Range("E51").Value = "  -1.43%  "
